$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the existing
# header style used by H1 (bold font, thin border, centered/top aligned).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I (I0) and J (IF) for rows 2-31.
$iVals = @(1,1,1,1,1,1,1,1,1,1,1,1,1,4,8,2,8,4,7,5,5,6,6,1,1,1,1,1,1,1)
$jVals = @(6,5,6,7,5,5,5,3,6,5,7,5,4,6,8,5,8,7,8,6,6,8,7,3,4,1,6,4,3,2)

for ($idx = 0; $idx -lt 30; $idx++) {
  $row = $idx + 2
  $ws.Cells.Item($row, 9).Value = $iVals[$idx]
  $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
